$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (constraint fraction k) across data rows 2-11
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary statistics table
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the new summary values: bold 12pt font, vertical-centered, row height taller
$summaryRange = $ws.Range("B14:B17")
$summaryRange.Font.Bold = $true
$summaryRange.Font.Size = 12
$summaryRange.VerticalAlignment = -4108

$ws.Range("A14:B17").EntireRow.RowHeight = 15.6

# Selection matching the authored workbook
$ws.Range("A14:B17").Select()

# Printer/page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
